$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("B9").Value = 90466

# Row 10
$ws.Range("A10").Value = 112178537
$ws.Range("B10").Value = 96720
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = "Knärot"
$ws.Range("G10").Value = "Goodyera repens"
$ws.Range("H10").Value = "(L.) R. Br."
$ws.Range("Q10").Value = 760382
$ws.Range("R10").Value = 7210147

# Row 11
$ws.Range("A11").Value = 112178540
$ws.Range("B11").Value = 90466
$ws.Range("Q11").Value = 760340
$ws.Range("R11").Value = 7210120

# Row 12
$ws.Range("A12").Value = 112178530
$ws.Range("B12").Value = 96720
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("Q12").Value = 760431
$ws.Range("R12").Value = 7210191

# Row 13
$ws.Range("A13").Value = 112178522
$ws.Range("B13").Value = 85434
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 3739
$ws.Range("F13").Value = "Persiljespindling"
$ws.Range("G13").Value = "Cortinarius sulfurinus"
$ws.Range("H13").Value = "Quél."
$ws.Range("Q13").Value = 760108
$ws.Range("R13").Value = 7210439

# Row 14
$ws.Range("A14").Value = 112178516
$ws.Range("B14").Value = 89317
$ws.Range("E14").Value = 3215
$ws.Range("F14").Value = "Rödgul trumpetsvamp"
$ws.Range("G14").Value = "Craterellus lutescens"
$ws.Range("H14").Value = "(Fr.) Fr."
$ws.Range("Q14").Value = 760126
$ws.Range("R14").Value = 7210471

# Row 15
$ws.Range("A15").Value = 112178514
$ws.Range("B15").Value = 102166
$ws.Range("E15").Value = 222412
$ws.Range("F15").Value = "Tibast"
$ws.Range("G15").Value = "Daphne mezereum"
$ws.Range("H15").Value = "L."
$ws.Range("Q15").Value = 760068
$ws.Range("R15").Value = 7210453

# Row 16
$ws.Range("A16").Value = 112178529
$ws.Range("B16").Value = 90466
$ws.Range("Q16").Value = 760450
$ws.Range("R16").Value = 7210211

# Row 17
$ws.Range("A17").Value = 112178521
$ws.Range("B17").Value = 90466
$ws.Range("E17").Value = 4769
$ws.Range("F17").Value = "Svavelriska"
$ws.Range("G17").Value = "Lactarius scrobiculatus"
$ws.Range("H17").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q17").Value = 760097
$ws.Range("R17").Value = 7210441

# Row 18
$ws.Range("A18").Value = 112178519
$ws.Range("B18").Value = 85386
$ws.Range("E18").Value = 1988
$ws.Range("F18").Value = "Kryddspindling"
$ws.Range("G18").Value = "Cortinarius percomis"
$ws.Range("H18").Value = "Fr."
$ws.Range("Q18").Value = 760104
$ws.Range("R18").Value = 7210466

# Row 19
$ws.Range("A19").Value = 112178532
$ws.Range("B19").Value = 90466
$ws.Range("Q19").Value = 760411
$ws.Range("R19").Value = 7210179

# Row 20
$ws.Range("B20").Value = 85434

# Row 21
$ws.Range("A21").Value = 112178526
$ws.Range("B21").Value = 90800
$ws.Range("E21").Value = 4364
$ws.Range("F21").Value = "Dropptaggsvamp"
$ws.Range("G21").Value = "Hydnellum ferrugineum"
$ws.Range("H21").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q21").Value = 760256
$ws.Range("R21").Value = 7210384

# Row 22
$ws.Range("A22").Value = 112178528
$ws.Range("B22").Value = 90466
$ws.Range("D22").Value = "LC"
$ws.Range("E22").Value = 4769
$ws.Range("F22").Value = "Svavelriska"
$ws.Range("G22").Value = "Lactarius scrobiculatus"
$ws.Range("H22").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q22").Value = 760519
$ws.Range("R22").Value = 7210363

# Row 23
$ws.Range("A23").Value = 112178515
$ws.Range("B23").Value = 90466
$ws.Range("E23").Value = 4769
$ws.Range("F23").Value = "Svavelriska"
$ws.Range("G23").Value = "Lactarius scrobiculatus"
$ws.Range("H23").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q23").Value = 760089
$ws.Range("R23").Value = 7210467

# Row 24
$ws.Range("A24").Value = 112178531
$ws.Range("B24").Value = 96720
$ws.Range("D24").Value = "VU"
$ws.Range("E24").Value = 220787
$ws.Range("F24").Value = "Knärot"
$ws.Range("G24").Value = "Goodyera repens"
$ws.Range("H24").Value = "(L.) R. Br."
$ws.Range("Q24").Value = 760437
$ws.Range("R24").Value = 7210197

# Row 25
$ws.Range("A25").Value = 112178517
$ws.Range("B25").Value = 90466
$ws.Range("E25").Value = 4769
$ws.Range("F25").Value = "Svavelriska"
$ws.Range("G25").Value = "Lactarius scrobiculatus"
$ws.Range("H25").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q25").Value = 760128
$ws.Range("R25").Value = 7210459

# Row 26
$ws.Range("A26").Value = 112178539
$ws.Range("B26").Value = 90466
$ws.Range("Q26").Value = 760354
$ws.Range("R26").Value = 7210135

# Row 27
$ws.Range("A27").Value = 112178538
$ws.Range("B27").Value = 98872
$ws.Range("D27").Value = "LC"
$ws.Range("E27").Value = 222771
$ws.Range("F27").Value = "Svart trolldruva"
$ws.Range("G27").Value = "Actaea spicata"
$ws.Range("H27").Value = "L."
$ws.Range("Q27").Value = 760363
$ws.Range("R27").Value = 7210127

# Row 28
$ws.Range("A28").Value = 112178535
$ws.Range("B28").Value = 102166
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 222412
$ws.Range("F28").Value = "Tibast"
$ws.Range("G28").Value = "Daphne mezereum"
$ws.Range("H28").Value = "L."
$ws.Range("Q28").Value = 760389
$ws.Range("R28").Value = 7210155

# Row 29
$ws.Range("A29").Value = 112178520
$ws.Range("B29").Value = 96753
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 219874
$ws.Range("F29").Value = "Nattviol"
$ws.Range("G29").Value = "Platanthera bifolia"
$ws.Range("H29").Value = "(L.) Rich."
$ws.Range("Q29").Value = 760092
$ws.Range("R29").Value = 7210449

# Row 30
$ws.Range("A30").Value = 112178524
$ws.Range("B30").Value = 90812
$ws.Range("E30").Value = 4366
$ws.Range("F30").Value = "Skarp dropptaggsvamp"
$ws.Range("G30").Value = "Hydnellum peckii"
$ws.Range("H30").Value = "Banker"
$ws.Range("Q30").Value = 760203
$ws.Range("R30").Value = 7210420
